$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (s="1") from D1 onto the two new header cells E1 and F1,
# then set the header labels. D1's label shifts from "message_content" to
# "message_content_z", E1 becomes "message_content_a", F1 becomes "message_content".
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)

$ws.Range("D1").Value = "message_content_z"
$ws.Range("E1").Value = "message_content_a"
$ws.Range("F1").Value = "message_content"

# Per-row data: row number, column D value (numeric), column E value (numeric),
# column F sentiment label (text).
$data = @(
  @(2, 0.8, 0.8, "Positive"),
  @(3, 0.2, 0.5, "Neutral"),
  @(4, 0.2, 0.5, "Neutral"),
  @(5, 0.1, 0.5, "Neutral"),
  @(6, -0.8, -0.5, "Negative"),
  @(7, -0.5, 0.2, "Neutral"),
  @(8, 0.3, 0.5, "Neutral"),
  @(9, 0, 0, "Neutral"),
  @(10, 0.8, 1, "Positive"),
  @(11, 0.2, 0.5, "Neutral"),
  @(12, 0.8, 0.8, "Positive"),
  @(13, 0.5, 0.8, "Positive"),
  @(14, 0.5, 0.5, "Positive"),
  @(15, 0.1, 0.2, "Neutral"),
  @(16, -0.3, 0, "Neutral"),
  @(17, 0.5, 0.5, "Positive"),
  @(18, 0.7, 0.7, "Positive"),
  @(19, 0.1, 0.5, "Neutral"),
  @(20, 0.6, 0.8, "Positive"),
  @(21, 0, 0, "Neutral"),
  @(22, 0, 0.6, "Neutral"),
  @(23, 0.1, 0.5, "Neutral"),
  @(24, 0.5, 0.5, "Neutral"),
  @(25, -1, -1, "Negative"),
  @(26, 0, 0.5, "Neutral"),
  @(27, -1, 0, "Neutral"),
  @(28, 0.7, 0.7, "Positive"),
  @(29, 0.5, 0.7, "Positive"),
  @(30, 0.7, 0.8, "Positive"),
  @(31, 0.5, 0.5, "Positive"),
  @(32, 0.3, 0.5, "Neutral"),
  @(33, 0.1, 0.5, "Neutral"),
  @(34, 0.2, 0.5, "Neutral"),
  @(35, 0.6, 0.8, "Positive"),
  @(36, -0.5, -0.5, "Neutral"),
  @(37, 0.2, 0.5, "Neutral"),
  @(38, 0.2, 0.5, "Neutral"),
  @(39, 0.3, 0.5, "Neutral"),
  @(40, -1, 0, "Neutral"),
  @(41, -1, -1, "Negative"),
  @(42, 0.5, 0.7, "Positive"),
  @(43, 0, 0, "Neutral"),
  @(44, -1, 0, "Neutral"),
  @(45, 0, 0.5, "Neutral"),
  @(46, 0.5, 0.5, "Positive"),
  @(47, -0.5, 0.1, "Neutral"),
  @(48, 0.5, 0.5, "Neutral"),
  @(49, 0.5, 0.6, "Positive"),
  @(50, -1, 0.5, "neutral"),
  @(51, 0.2, 0.5, "Neutral"),
  @(52, 0.2, 0.5, "Neutral"),
  @(53, 0, 0.5, "Neutral"),
  @(54, 0.6, 0.6, "Positive"),
  @(55, 0.2, 0.5, "Neutral"),
  @(56, -0.5, 0, "Negative"),
  @(57, 0.5, 0.5, "Positive"),
  @(58, -1, 0, "Neutral"),
  @(59, 0.7, 0.5, "Positive"),
  @(60, 0.5, 0.6, "Positive"),
  @(61, 0.2, 0.5, "Neutral"),
  @(62, -1, 0, "Neutral"),
  @(63, 0.5, 0.5, "Neutral"),
  @(64, 0.8, 0.8, "Positive"),
  @(65, 0.2, 0.3, "Neutral"),
  @(66, 0.2, 0.8, "Positive"),
  @(67, 0.8, 0.8, "Positive"),
  @(68, 0.8, 0.9, "Positive"),
  @(69, 0.6, 0.6, "Neutral"),
  @(70, -0.5, -0.5, "Negative"),
  @(71, -0.2, 0, "Neutral"),
  @(72, 0.2, 0.5, "Positive"),
  @(73, 0.7, 0.8, "Positive"),
  @(74, 0.5, 0.6, "Positive"),
  @(75, 0.5, 0.6, "Neutral"),
  @(76, 0.8, 0.7, "Positive"),
  @(77, 0.7, 0.7, "Neutral"),
  @(78, -0.5, 0.2, "Neutral"),
  @(79, -0.5, 0, "Neutral")
)

foreach ($row in $data) {
  $r = $row[0]
  $dVal = $row[1]
  $eVal = $row[2]
  $fVal = $row[3]

  $ws.Cells.Item($r, 4).Value = $dVal
  $ws.Cells.Item($r, 5).Value = $eVal
  $ws.Cells.Item($r, 6).Value = $fVal
}

Write-Output "applied sentiment columns E and F for rows 2-79"
